$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5, 13, 6, 7),
    @(4, 7, 5, 13),
    @(4, 4, 6, 16),
    @(4, 12, 3, 8),
    @(4, 8, 2, 12),
    @(5, 8, 3, 12),
    @(3, 13, 2, 7),
    @(2, 15, 5, 5),
    @(9, 6, 7, 14),
    @(3, 15, 5, 5),
    @(4, 17, 5, 3),
    @(4, 13, 3, 7),
    @(6, 16, 4, 4),
    @(5, 6, 7, 14),
    @(5, 14, 6, 6),
    @(4, 6, 2, 14)
)

$startRow = 1019
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}

$ws.Application.ActiveWindow.ScrollRow = 1013
$ws.Range("A1035").Select()

$wb.Save()
